# Apply the commit's changes:
#  - sheet "full_signals - without decay": update view selection
#  - sheet "full_signals - with decay": remove the "warm_steps(20%)" column,
#    rename the remaining decay column, update model-1 row values, and
#    update the view selection

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("full_signals - without decay")
$ws2 = $wb.Worksheets.Item("full_signals - with decay")

# --- Sheet 1: "full_signals - without decay" -------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A28:XFD32").Select() | Out-Null

# --- Sheet 2: "full_signals - with decay" -----------------------------------
$ws2.Activate() | Out-Null

# Remove the old "warm_steps(20%)" column; everything to the right shifts left.
$ws2.Range("G:G").Delete() | Out-Null

# The former "decay_steps (80%)" column (now column G) becomes just "decay_steps".
$ws2.Range("G5").Value2 = "decay_steps"

# Update the model-1 results row with the new Cosine Decay run data.
$ws2.Range("B6").Value2 = "modelo 1"
$ws2.Range("E6").Value2 = 500
$ws2.Range("K6").Value2 = 0.0714
$ws2.Range("L6").Value2 = 0.1849
$ws2.Range("M6").Value2 = 0.0744
$ws2.Range("N6").Value2 = 0.164
$ws2.Range("P6").Value2 = 116

$ws2.Range("L20").Select() | Out-Null
